$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.578.86'
$ws.Range('E2').Value = '  -1.01%  '
$ws.Range('D3').Value = '1.598.03'
$ws.Range('E3').Value = '  -1.82%  '
$ws.Range('E4').Value = '  +0.38%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '208.34'
$ws.Range('E5').Value = '  -1.06%  '
$ws.Range('E6').Value = '  -3.35%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '22.32'
$ws.Range('E8').Value = '  -4.42%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.253'
$ws.Range('E9').Value = '  -1.71%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0593'
$ws.Range('E10').Value = '  -3.14%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0867'
$ws.Range('E11').Value = '  -1.42%  '
$ws.Range('D12').Value = '1.824.20'
$ws.Range('E12').Value = '  -1.82%  '
$ws.Range('D13').Value = '1.595.00'
$ws.Range('E13').Value = '  -1.97%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.88'
$ws.Range('E14').Value = '  -3.69%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.541'
$ws.Range('E15').Value = '  -3.58%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '63.48'
$ws.Range('E16').Value = '  -2.83%  '
$ws.Range('D17').Value = '27.544.49'
$ws.Range('E17').Value = '  -1.08%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '219.37'
$ws.Range('E18').Value = '  -4.53%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.41'
$ws.Range('E19').Value = '  -3.11%  '
$ws.Range('D20').Value = '0.0₃0696'
$ws.Range('E20').Value = '  -3.36%  '
$ws.Range('E21').Value = '  +0.47%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.22'
$ws.Range('E22').Value = '  -2.34%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.71'
$ws.Range('E23').Value = '  -3.70%  '
$ws.Range('E24').Value = '  -2.42%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '154.56'
$ws.Range('E25').Value = '  +0.09%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.77'
$ws.Range('E26').Value = '  -1.93%  '
$ws.Range('E27').Value = '  +0.39%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.09'
$ws.Range('E28').Value = '  -2.77%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.107'
$ws.Range('E29').Value = '  -3.90%  '
$ws.Range('E30').Value = '  -0.90%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0470'
$ws.Range('E31').Value = '  -2.41%  '
$ws.Range('E32').Value = '  -4.32%  '
$ws.Range('D33').Value = '1.365.56'
$ws.Range('E33').Value = '  -2.24%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.96'
$ws.Range('E35').Value = '  -2.67%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.972'
$ws.Range('E36').Value = '  -4.04%  '
$ws.Range('E37').Value = '  -0.88%  '
$ws.Range('E38').Value = '  -2.30%  '
$ws.Range('E39').Value = '  -2.62%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.816'
$ws.Range('E40').Value = '  -3.95%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.980'
$ws.Range('E42').Value = '  -2.29%  '
$ws.Range('E43').Value = '  -0.95%  '
$ws.Range('B44').Value = 'RenderToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.78'
$ws.Range('E44').Value = '  -3.34%  '
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '64.13'
$ws.Range('E45').Value = '  -2.47%  '
$ws.Range('D46').Value = '1.733.65'
$ws.Range('E46').Value = '  -1.93%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '88.00'
$ws.Range('E48').Value = '  -0.04%  '
$ws.Range('B49').Value = 'Algorand'
$ws.Range('C49').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0974'
$ws.Range('E49').Value = '  -4.00%  '
$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').Value = '0.0₇0981'
$ws.Range('E50').Value = '  -3.73%  '
$ws.Range('E51').Value = '  -1.15%  '
